$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.140.94'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.418.11'
$ws.Range("E3").Value = '  -0.54%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '552.18'
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.10'
$ws.Range("E6").Value = '  -1.57%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +1.66%  '
$ws.Range("E9").Value = '  -1.36%  '
$ws.Range("E10").Value = '  -1.88%  '
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.848.81'
$ws.Range("E14").Value = '  -0.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.044.54'
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("E16").Value = '  -1.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.417.32'
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("E19").Value = '  +1.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '326.87'
$ws.Range("E20").Value = '  -1.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.75'
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.45'
$ws.Range("E23").Value = '  +0.38%  '
$ws.Range("E24").Value = '  +2.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.69'
$ws.Range("E25").Value = '  +0.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.996'
$ws.Range("E26").Value = '  -0.46%  '
$ws.Range("E27").Value = '  +2.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0772'
$ws.Range("E28").Value = '  -1.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.76'
$ws.Range("E29").Value = '  -1.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.96'
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("E31").Value = '  -3.93%  '
$ws.Range("E32").Value = '  +2.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.403'
$ws.Range("E33").Value = '  -4.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.55'
$ws.Range("E34").Value = '  -1.14%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  +0.49%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '329.55'
$ws.Range("E39").Value = '  +1.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.81'
$ws.Range("E41").Value = '  -2.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '145.14'
$ws.Range("E42").Value = '  +3.21%  '
$ws.Range("E43").Value = '  -1.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.02'
$ws.Range("E44").Value = '  +2.13%  '
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("E46").Value = '  -2.13%  '
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("E48").Value = '  -1.95%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '11.05'
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("E50").Value = '  -3.60%  '
$ws.Range("E51").Value = '  -1.03%  '
